$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$src = $ws1.Range("A2")
$src.Copy()
$dst = $ws1.Range("A4")
$dst.PasteSpecial(-4122)
